$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Participant name placeholder ("NOM / PRENOM" run keeps its own run,
#    neighbouring runs use a different font so a plain replace cannot bleed
#    into them).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("NOM / PRENOM", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{participant_prenom}} {{participant_nom}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Formation title placeholder (single-run paragraph, keep trailing space).
# ---------------------------------------------------------------------------
$titreOld = "Renouvellement CACES R386 CAT" + [char]0x00C9 + "GORIE 3B "
$d.Content.Find.Execute($titreOld, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{formation_titre}} ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Formation duration placeholder (single-run paragraph).
# ---------------------------------------------------------------------------
$rngDuree = $d.Content.Duplicate
$rngDuree.Find.Execute("14", $true, $false, $false, $false, $false, `
                        $true, 1, $false, $null, 0) | Out-Null
$rngDuree.Text = "{{formation_duree}}"

# ---------------------------------------------------------------------------
# 4) Start date placeholder (single-run paragraph).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("15/11/2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{date_debut}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) End date placeholder (single-run paragraph).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("22/11/2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{date_fin}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Start time placeholder. This run lives in a 3-run paragraph
#    ("8H30 -" / "1" / "6H30") where all three runs share identical
#    formatting. A direct replace of the first run causes this host to
#    coalesce it with its following siblings (they'd lose their own
#    <w:r> boundaries). To keep "1" and "6H30" as their own untouched
#    runs, we briefly nudge the very next run's formatting (Bold on/off)
#    so it no longer matches the edited run, edit the text, then restore
#    the original formatting - which does not re-trigger a merge.
# ---------------------------------------------------------------------------
$rngGuard = $d.Content.Duplicate
$rngGuard.Find.Execute("8H30 -", $true, $false, $false, $false, $false, `
                        $true, 1, $false, $null, 0) | Out-Null
$guardStart = $rngGuard.End
$guardRange = $d.Range($guardStart, $guardStart + 1)
$guardRange.Font.Bold = 1

$rngHoraire = $d.Content.Duplicate
$rngHoraire.Find.Execute("8H30 -", $true, $false, $false, $false, $false, `
                          $true, 1, $false, $null, 0) | Out-Null
$rngHoraire.Text = "{{horaire_debut}} -"

$rngRestore = $d.Content.Duplicate
$rngRestore.Find.Execute("{{horaire_debut}} -", $true, $false, $false, $false, $false, `
                          $true, 1, $false, $null, 0) | Out-Null
$restorePos = $rngRestore.End
$d.Range($restorePos, $restorePos + 1).Font.Bold = 0

# ---------------------------------------------------------------------------
# 7) Venue name placeholder (single-run paragraph).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Entreprise MAHEY", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{lieu}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Street-address line: the run's text is removed entirely (the <w:t>
#    element disappears) while the run (and its <w:rPr>) and the
#    paragraph's own <w:pPr> stay in place. A plain text delete prunes the
#    now-empty run altogether, so rebuild the paragraph explicitly via
#    InsertXML, preserving its original paragraph/run formatting.
# ---------------------------------------------------------------------------
$rngAdresse = $d.Content.Duplicate
$rngAdresse.Find.Execute("5, Impasse Grand Jardin", $true, $false, $false, $false, $false, `
                          $true, 1, $false, $null, 0) | Out-Null
$adresseXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:pPr><w:ind w:right="1"/><w:rPr><w:rFonts w:ascii="ITC Avant Garde Gothic Book" w:hAnsi="ITC Avant Garde Gothic Book" w:cs="Calibri"/><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="fr-FR"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="ITC Avant Garde Gothic Book" w:hAnsi="ITC Avant Garde Gothic Book" w:cs="Calibri"/><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="fr-FR"/></w:rPr></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$rngAdresse.InsertXML($adresseXml) | Out-Null

# ---------------------------------------------------------------------------
# 9) Postal code / city placeholder (single-run paragraph).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("35400 Saint Malo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{organisme_code_postal}} {{organisme_ville}}", 2) | Out-Null
